# Auto-generated edit script applying Asura_Profits.xlsx market-data refresh
# Updates cached price/profit figures on ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1607.7778
$ws.Range("J40").Value = 1762.5
$ws.Range("L40").Value = 1762.5
$ws.Range("N40").Value = -2112.5
$ws.Range("H41").Value = 320
$ws.Range("I41").Value = 358.41666
$ws.Range("J41").Value = 268.77777
$ws.Range("K41").Value = 358.41666
$ws.Range("L41").Value = 268.77777
$ws.Range("M41").Value = 81.58334000000002
$ws.Range("N41").Value = -1148.77777
$ws.Range("H64").Value = 3303
$ws.Range("I64").Value = 3148.4666
$ws.Range("J64").Value = 3425
$ws.Range("K64").Value = 3148.4666
$ws.Range("L64").Value = 3425
$ws.Range("M64").Value = -2900.4666
$ws.Range("N64").Value = -3921
$ws.Range("H67").Value = 3303
$ws.Range("I67").Value = 3148.4666
$ws.Range("J67").Value = 3425
$ws.Range("K67").Value = 3148.4666
$ws.Range("L67").Value = 3425
$ws.Range("M67").Value = -2290.4666
$ws.Range("N67").Value = -5141
$ws.Range("H74").Value = 3834.7368
$ws.Range("I74").Value = 3512.5
$ws.Range("J74").Value = 3983.4614
$ws.Range("K74").Value = 3512.5
$ws.Range("L74").Value = 3983.4614
$ws.Range("M74").Value = -2576.5
$ws.Range("N74").Value = -5855.4614
$ws.Range("H77").Value = 3834.7368
$ws.Range("I77").Value = 3512.5
$ws.Range("J77").Value = 3983.4614
$ws.Range("K77").Value = 17562.5
$ws.Range("L77").Value = 19917.307
$ws.Range("M77").Value = -12882.5
$ws.Range("N77").Value = -29277.307
$ws.Range("H140").Value = 82406.25
$ws.Range("J140").Value = 88821.42999999999
$ws.Range("L140").Value = 88821.42999999999
$ws.Range("N140").Value = -99181.42999999999

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9888.522000000001
$ws.Range("I32").Value = 9640.843999999999
$ws.Range("K32").Value = 9640.843999999999
$ws.Range("M32").Value = -9353.843999999999
$ws.Range("H88").Value = 1990.1111
$ws.Range("I88").Value = 1300
$ws.Range("J88").Value = 2335.1667
$ws.Range("K88").Value = 1300
$ws.Range("L88").Value = 2335.1667
$ws.Range("M88").Value = -894
$ws.Range("N88").Value = -3147.1667
$ws.Range("H91").Value = 1990.1111
$ws.Range("I91").Value = 1300
$ws.Range("J91").Value = 2335.1667
$ws.Range("K91").Value = 1300
$ws.Range("L91").Value = 2335.1667
$ws.Range("M91").Value = 104
$ws.Range("N91").Value = -5143.1667
$ws.Range("H123").Value = 48429
$ws.Range("J123").Value = 48429
$ws.Range("L123").Value = 48429
$ws.Range("N123").Value = -58229

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H102").Value = 168418.67
$ws.Range("I102").Value = 168418.67
$ws.Range("K102").Value = 168418.67
$ws.Range("M102").Value = -165173.67
$ws.Range("H105").Value = 1307
$ws.Range("I105").Value = 1307
$ws.Range("K105").Value = 1307
$ws.Range("M105").Value = 440
$ws.Range("H134").Value = 2456.7812
$ws.Range("I134").Value = 2430.8845
$ws.Range("J134").Value = 2569
$ws.Range("K134").Value = 7292.6535
$ws.Range("L134").Value = 7707
$ws.Range("M134").Value = -4757.6535
$ws.Range("N134").Value = -12777

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23258684
$ws.Range("J31").Value = 4392.0713
$ws.Range("L31").Value = 4392.0713
$ws.Range("N31").Value = -4982.0713
$ws.Range("H34").Value = 23258684
$ws.Range("J34").Value = 4392.0713
$ws.Range("L34").Value = 4392.0713
$ws.Range("N34").Value = -4796.0713
$ws.Range("H62").Value = 49662.727
$ws.Range("I62").Value = 73827.14
$ws.Range("J62").Value = 7375
$ws.Range("K62").Value = 73827.14
$ws.Range("L62").Value = 7375
$ws.Range("M62").Value = -73203.14
$ws.Range("N62").Value = -8623
$ws.Range("H65").Value = 49662.727
$ws.Range("I65").Value = 73827.14
$ws.Range("J65").Value = 7375
$ws.Range("K65").Value = 369135.7
$ws.Range("L65").Value = 36875
$ws.Range("M65").Value = -366015.7
$ws.Range("N65").Value = -43115
$ws.Range("H132").Value = 3150.6875
$ws.Range("I132").Value = 2701
$ws.Range("J132").Value = 4499.75
$ws.Range("K132").Value = 8103
$ws.Range("L132").Value = 13499.25
$ws.Range("M132").Value = -5573
$ws.Range("N132").Value = -18559.25

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 600
$ws.Range("I44").Value = 60
$ws.Range("J44").Value = 735
$ws.Range("K44").Value = 180
$ws.Range("L44").Value = 2205
$ws.Range("M44").Value = 218
$ws.Range("N44").Value = -3001

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H58").Value = 10000
$ws.Range("I58").Value = 10000
$ws.Range("K58").Value = 10000
$ws.Range("M58").Value = -9723
$ws.Range("H70").Value = 6172.5557
$ws.Range("I70").Value = 6024.75
$ws.Range("J70").Value = 6290.8
$ws.Range("K70").Value = 6024.75
$ws.Range("L70").Value = 6290.8
$ws.Range("M70").Value = -5754.75
$ws.Range("N70").Value = -6830.8
$ws.Range("H73").Value = 6172.5557
$ws.Range("I73").Value = 6024.75
$ws.Range("J73").Value = 6290.8
$ws.Range("K73").Value = 6024.75
$ws.Range("L73").Value = 6290.8
$ws.Range("M73").Value = -5088.75
$ws.Range("N73").Value = -8162.8
$ws.Range("H80").Value = 3257.4285
$ws.Range("I80").Value = 3456.111
$ws.Range("K80").Value = 3456.111
$ws.Range("M80").Value = -2458.111
$ws.Range("H83").Value = 3257.4285
$ws.Range("I83").Value = 3456.111
$ws.Range("K83").Value = 17280.555
$ws.Range("M83").Value = -12288.555

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3330.5454
$ws.Range("I7").Value = 2600.1428
$ws.Range("K7").Value = 2600.1428
$ws.Range("M7").Value = -2488.1428
$ws.Range("H82").Value = 2779.4167
$ws.Range("I82").Value = 1821.4286
$ws.Range("J82").Value = 4120.6
$ws.Range("K82").Value = 1821.4286
$ws.Range("L82").Value = 4120.6
$ws.Range("M82").Value = -1460.4286
$ws.Range("N82").Value = -4842.6
$ws.Range("H85").Value = 2779.4167
$ws.Range("I85").Value = 1821.4286
$ws.Range("J85").Value = 4120.6
$ws.Range("K85").Value = 1821.4286
$ws.Range("L85").Value = 4120.6
$ws.Range("M85").Value = -573.4286
$ws.Range("N85").Value = -6616.6
$ws.Range("H93").Value = 1801.7858
$ws.Range("I93").Value = 1748.0769
$ws.Range("J93").Value = 2500
$ws.Range("K93").Value = 1748.0769
$ws.Range("L93").Value = 2500
$ws.Range("M93").Value = -500.0769
$ws.Range("N93").Value = -4996
$ws.Range("H126").Value = 3330.5454
$ws.Range("I126").Value = 2600.1428
$ws.Range("K126").Value = 7800.428400000001
$ws.Range("M126").Value = -5330.428400000001

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 34500
$ws.Range("J54").Value = 34500
$ws.Range("L54").Value = 34500
$ws.Range("N54").Value = -35540
